$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 39, shifting the existing rows 39-76 down to 40-77.
$ws.Rows.Item(39).EntireRow.Insert()

# Populate the newly inserted row 39 with the new weekly price-report record.
$ws.Range("A39").Value = 10
$ws.Range("B39").Value = "Vega Modelo de Temuco"
$ws.Range("C39").Value = "La Araucanía"
$ws.Range("D39").Value = 44827
$ws.Range("D39").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E39").Value = 9
$ws.Range("F39").Value = 100112026
$ws.Range("G39").Value = "Haba"
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 50
$ws.Range("K39").Value = 10000
$ws.Range("L39").Value = 10000
$ws.Range("M39").Value = 10000
$ws.Range("N39").Value = "$/saco 25 kilos"
$ws.Range("O39").Value = "Provincia de Limarí"
$ws.Range("P39").Value = 400
$ws.Range("Q39").Value = 25
$ws.Range("R39").Value = "Hortaliza"
